$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Nitrate - rename only
$ws.Range("A2").Value = "Nitrate_c_Night_sp_exchange"

# Row 3: L-Valine - rename + update C3, D3
$ws.Range("A3").Value = "L-Valine_Night_sp_exchange"
$ws.Range("C3").Value = -0.02090000000000001
$ws.Range("D3").Value = 0.02090000000001012

# Row 4: L-Proline - rename + update C4, D4
$ws.Range("A4").Value = "L-Proline_Night_sp_exchange"
$ws.Range("C4").Value = -0.2676579583333158
$ws.Range("D4").Value = 3.966457458332928

# Row 5: L-Alanine_c - rename + update C5
$ws.Range("A5").Value = "L-Alanine_c_Night_sp_exchange"
$ws.Range("C5").Value = -9.637505166672561

# Row 6: D-Fructose - rename + update C6, D6
$ws.Range("A6").Value = "D-Fructose_Night_sp_exchange"
$ws.Range("C6").Value = -0.5590199999950832
$ws.Range("D6").Value = 0.5590200000001522

# Row 7: Citrate_c - rename + update C7, D7
$ws.Range("A7").Value = "Citrate_c_Night_sp_exchange"
$ws.Range("C7").Value = -9.747673833332518
$ws.Range("D7").Value = -1.348724666666926
